$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1416.5264
$ws.Range("I28").Value = 1416.5264
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1416.5264
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -931.5264
$ws.Range("N28").ClearContents()

$ws.Range("H53").Value = 230.29411
$ws.Range("I53").Value = 161.5
$ws.Range("J53").Value = 328.57144
$ws.Range("K53").Value = 161.5
$ws.Range("L53").Value = 328.57144
$ws.Range("M53").Value = 475.5

$ws.Range("H62").Value = 1919
$ws.Range("I62").Value = 1919
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1919
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1295
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 1919
$ws.Range("I65").Value = 1919
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9595
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6475
$ws.Range("N65").ClearContents()

$ws.Range("H76").Value = 3678.8333
$ws.Range("I76").Value = 3797
$ws.Range("J76").Value = 3088
$ws.Range("K76").Value = 3797
$ws.Range("L76").Value = 3088
$ws.Range("M76").Value = -3482
$ws.Range("N76").Value = -3718

$ws.Range("H79").Value = 3678.8333
$ws.Range("I79").Value = 3797
$ws.Range("J79").Value = 3088
$ws.Range("K79").Value = 3797
$ws.Range("L79").Value = 3088
$ws.Range("M79").Value = -2705
$ws.Range("N79").Value = -5272

$ws.Range("H106").Value = 26371.143
$ws.Range("I106").Value = 30483
$ws.Range("J106").Value = 1700
$ws.Range("K106").Value = 30483
$ws.Range("L106").Value = 1700
$ws.Range("M106").Value = -29852

$ws.Range("H125").Value = 88239140
$ws.Range("I125").Value = 100003416
$ws.Range("J125").Value = 71433020
$ws.Range("K125").Value = 900030744
$ws.Range("L125").Value = 642897180
$ws.Range("M125").Value = -900028284

$ws.Range("H127").Value = 1839.3846
$ws.Range("I127").Value = 2187.8
$ws.Range("J127").Value = 1621.625
$ws.Range("K127").Value = 6563.400000000001
$ws.Range("L127").Value = 4864.875
$ws.Range("M127").Value = -1603.400000000001

$ws.Range("H129").Value = 2559.261
$ws.Range("I129").Value = 1046
$ws.Range("J129").Value = 2877.842
$ws.Range("K129").Value = 3138
$ws.Range("L129").Value = 8633.526
$ws.Range("M129").Value = 1862
$ws.Range("N129").Value = -18633.526

$ws.Range("H131").Value = 20927.8
$ws.Range("I131").Value = 20927.8
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 62783.39999999999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -57743.39999999999

$ws.Range("H137").Value = 1797.2273
$ws.Range("I137").Value = 1267
$ws.Range("J137").Value = 3600
$ws.Range("K137").Value = 3801
$ws.Range("L137").Value = 10800
$ws.Range("M137").Value = -1251
$ws.Range("N137").Value = -15900

$ws.Range("H141").Value = 3582.6667
$ws.Range("I141").Value = 2874.25
$ws.Range("J141").Value = 4999.5
$ws.Range("K141").Value = 8622.75
$ws.Range("L141").Value = 14998.5
$ws.Range("M141").Value = -3442.75
$ws.Range("N141").Value = -25358.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2898.625
$ws.Range("I16").Value = 3539.4
$ws.Range("J16").Value = 1830.6666
$ws.Range("K16").Value = 3539.4
$ws.Range("L16").Value = 1830.6666
$ws.Range("M16").Value = -3252.4
$ws.Range("N16").Value = -2404.6666

$ws.Range("H22").Value = 675.125
$ws.Range("I22").Value = 608.8
$ws.Range("J22").Value = 785.6667
$ws.Range("K22").Value = 608.8
$ws.Range("L22").Value = 785.6667
$ws.Range("M22").Value = -258.8
$ws.Range("N22").Value = -1485.6667

$ws.Range("H107").Value = 964.13336
$ws.Range("I107").Value = 891.63635
$ws.Range("J107").Value = 1163.5
$ws.Range("K107").Value = 891.63635
$ws.Range("L107").Value = 1163.5
$ws.Range("M107").Value = 1028.36365
$ws.Range("N107").Value = -5003.5

$ws.Range("H113").Value = 2898.625
$ws.Range("I113").Value = 3539.4
$ws.Range("J113").Value = 1830.6666
$ws.Range("K113").Value = 3539.4
$ws.Range("L113").Value = 1830.6666
$ws.Range("M113").Value = -1369.4
$ws.Range("N113").Value = -6170.6666

$ws.Range("H134").Value = 2464.6
$ws.Range("I134").Value = 2441
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 7323
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -4788

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H141").Value = 103193.93
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 103193.93
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 103193.93
$ws.Range("N141").Value = -113553.93

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1610.4286
$ws.Range("I68").Value = 400
$ws.Range("J68").Value = 1812.1666
$ws.Range("K68").Value = 1200
$ws.Range("L68").Value = 5436.4998
$ws.Range("M68").Value = -389
$ws.Range("N68").Value = -7058.4998

$ws.Range("H71").Value = 1610.4286
$ws.Range("I71").Value = 400
$ws.Range("J71").Value = 1812.1666
$ws.Range("K71").Value = 3600
$ws.Range("L71").Value = 16309.4994
$ws.Range("M71").Value = 456
$ws.Range("N71").Value = -24421.4994

$ws.Range("H80").Value = 9500
$ws.Range("I80").Value = 9500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 28500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -27564

$ws.Range("H83").Value = 9500
$ws.Range("I83").Value = 9500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 85500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -80820

$ws.Range("H92").Value = 600
$ws.Range("I92").Value = 700
$ws.Range("J92").Value = 550
$ws.Range("K92").Value = 2100
$ws.Range("L92").Value = 1650
$ws.Range("M92").Value = -852
$ws.Range("N92").Value = -4146

$ws.Range("H121").Value = 9430.388999999999
$ws.Range("I121").Value = 38453.332
$ws.Range("J121").Value = 3625.8
$ws.Range("K121").Value = 115359.996
$ws.Range("L121").Value = 10877.4
$ws.Range("M121").Value = -114049.996

$ws.Range("H129").Value = 1003798.2
$ws.Range("I129").Value = 2816.6667
$ws.Range("J129").Value = 1432790.2
$ws.Range("K129").Value = 8450.000100000001
$ws.Range("L129").Value = 4298370.6
$ws.Range("M129").Value = -3450.000100000001
$ws.Range("N129").Value = -4308370.6

$ws.Range("H131").Value = 456788.78
$ws.Range("I131").Value = 1167.2858
$ws.Range("J131").Value = 669412.1
$ws.Range("K131").Value = 3501.8574
$ws.Range("L131").Value = 2008236.3
$ws.Range("M131").Value = 1538.1426
$ws.Range("N131").Value = -2018316.3

$ws.Range("H139").Value = 572.5714
$ws.Range("I139").Value = 402
$ws.Range("J139").Value = 999
$ws.Range("K139").Value = 1206
$ws.Range("L139").Value = 2997
$ws.Range("M139").Value = 3934
$ws.Range("N139").Value = -13277

$ws.Range("H140").Value = 12190.4
$ws.Range("I140").Value = 3554.3333
$ws.Range("J140").Value = 14349.417
$ws.Range("K140").Value = 10662.9999
$ws.Range("L140").Value = 43048.251
$ws.Range("M140").Value = -5482.999899999999
$ws.Range("N140").Value = -53408.251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8530
$ws.Range("I7").Value = 2450
$ws.Range("J7").Value = 9635.454
$ws.Range("K7").Value = 2450
$ws.Range("L7").Value = 9635.454
$ws.Range("M7").Value = -2338
$ws.Range("N7").Value = -9859.454

$ws.Range("H61").Value = 2989.8333
$ws.Range("I61").Value = 2987.8
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2987.8
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2785.8

$ws.Range("H113").Value = 2989.8333
$ws.Range("I113").Value = 2987.8
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2987.8
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -817.8000000000002

$ws.Range("H126").Value = 8530
$ws.Range("I126").Value = 2450
$ws.Range("J126").Value = 9635.454
$ws.Range("K126").Value = 7350
$ws.Range("L126").Value = 28906.362
$ws.Range("M126").Value = -4880
$ws.Range("N126").Value = -33846.362

$ws.Range("H132").Value = 4929.8096
$ws.Range("I132").Value = 5030.0586
$ws.Range("J132").Value = 4503.75
$ws.Range("K132").Value = 15090.1758
$ws.Range("L132").Value = 13511.25
$ws.Range("M132").Value = -12560.1758

$ws.Range("H134").Value = 72000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 72000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 72000
$ws.Range("N134").Value = -82140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 31000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 31000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 31000
$ws.Range("N44").Value = -32108

$ws.Range("H113").Value = 399
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 399
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1197
$ws.Range("N113").Value = -5537

$ws.Range("H136").Value = 2973.353
$ws.Range("I136").Value = 3113.1667
$ws.Range("J136").Value = 2637.8
$ws.Range("K136").Value = 9339.500100000001
$ws.Range("L136").Value = 7913.400000000001
$ws.Range("M136").Value = -6789.500100000001
$ws.Range("N136").Value = -13013.4
